$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 516 - this shifts the existing rows 516:609 down to 517:610
# and grows the used range from A1:R609 to A1:R610.
$ws.Rows.Item(516).Insert()

# Populate the newly inserted row 516 with the new record's data.
$ws.Range("A516").Value2 = 5
$ws.Range("B516").Value2 = "Macroferia Regional de Talca"
$ws.Range("C516").Value2 = "Maule"
$ws.Range("D516").Value2 = 44995
$ws.Range("E516").Value2 = 7
$ws.Range("F516").Value2 = 100112043
$ws.Range("G516").Value2 = "Pepino ensalada"
$ws.Range("H516").Value2 = "Sin especificar"
$ws.Range("I516").Value2 = "Primera"
$ws.Range("J516").Value2 = 300
$ws.Range("K516").Value2 = 10000
$ws.Range("L516").Value2 = 10000
$ws.Range("M516").Value2 = 10000
$ws.Range("N516").Value2 = "`$/caja 80 unidades"
$ws.Range("O516").Value2 = "Región del Maule"
$ws.Range("P516").Value2 = 125
$ws.Range("Q516").Value2 = 80
$ws.Range("R516").Value2 = "Hortaliza"
